{"js": "// Update the date line and the 25 two-digit-division answers in the table.\n// Each table cell's new value is written directly (by row/column index) so\n// that cross-matching text between \"before\" and \"after\" values (e.g. one\n// cell's new answer equals another cell's old answer) can never cause a\n// stray double-replacement.\n\n// 1. Update the date paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].insertText(\"2024-05-15 Wednesday\", Word.InsertLocation.replace);\n\n// 2. Update the practice-problem table, cell by cell.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, col, newText) - only the four \"content\" rows (0, 4, 8, 12, 16) of\n// the 20-row table hold text; the rows between them are spacer rows.\nconst updates = [\n  [0, 0, \"63\u00f75=12, 3\"],\n  [0, 1, \"81\u00f76=13, 3\"],\n  [0, 2, \"42\u00f75=8, 2\"],\n  [0, 3, \"58\u00f75=11, 3\"],\n  [0, 4, \"31\u00f79=3, 4\"],\n\n  [4, 0, \"66\u00f77=9, 3\"],\n  [4, 1, \"71\u00f75=14, 1\"],\n  [4, 2, \"16\u00f77=2, 2\"],\n  [4, 3, \"88\u00f75=17, 3\"],\n  [4, 4, \"44\u00f77=6, 2\"],\n\n  [8, 0, \"26\u00f78=3, 2\"],\n  [8, 1, \"68\u00f75=13, 3\"],\n  [8, 2, \"42\u00f74=10, 2\"],\n  [8, 3, \"17\u00f75=3, 2\"],\n  [8, 4, \"79\u00f74=19, 3\"],\n\n  [12, 0, \"46\u00f73=15, 1\"],\n  [12, 1, \"87\u00f72=43, 1\"],\n  [12, 2, \"87\u00f74=21, 3\"],\n  [12, 3, \"46\u00f77=6, 4\"],\n  [12, 4, \"78\u00f72=39, 0\"],\n\n  [16, 0, \"50\u00f78=6, 2\"],\n  [16, 1, \"89\u00f76=14, 5\"],\n  [16, 2, \"46\u00f79=5, 1\"],\n  [16, 3, \"36\u00f77=5, 1\"],\n  [16, 4, \"57\u00f72=28, 1\"],\n];\n\nfor (const [row, col, text] of updates) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 two-digit-division answers in the table.\n# Each table cell is addressed directly by its (row, column) position (Word's\n# Table.Cell is 1-indexed) so that cross-matching text between \"before\" and\n# \"after\" values (one cell's new answer can equal another cell's old answer)\n# never causes an accidental double replacement.\n\n$d = $word.ActiveDocument\n\n# 1. Update the date paragraph (first paragraph in the document).\n$d.Paragraphs.Item(1).Range.Text = \"2024-05-15 Wednesday\"\n\n# 2. Update the practice-problem table, cell by cell.\n$t = $d.Tables.Item(1)\n\n# (row, col, newText) - Word table rows/cols are 1-indexed. Only the four\n# \"content\" rows (1, 5, 9, 13, 17) of the 20-row table hold text; the rows\n# between them are spacer rows.\n$updates = @(\n    @(1, 1, \"63\u00f75=12, 3\"),\n    @(1, 2, \"81\u00f76=13, 3\"),\n    @(1, 3, \"42\u00f75=8, 2\"),\n    @(1, 4, \"58\u00f75=11, 3\"),\n    @(1, 5, \"31\u00f79=3, 4\"),\n\n    @(5, 1, \"66\u00f77=9, 3\"),\n    @(5, 2, \"71\u00f75=14, 1\"),\n    @(5, 3, \"16\u00f77=2, 2\"),\n    @(5, 4, \"88\u00f75=17, 3\"),\n    @(5, 5, \"44\u00f77=6, 2\"),\n\n    @(9, 1, \"26\u00f78=3, 2\"),\n    @(9, 2, \"68\u00f75=13, 3\"),\n    @(9, 3, \"42\u00f74=10, 2\"),\n    @(9, 4, \"17\u00f75=3, 2\"),\n    @(9, 5, \"79\u00f74=19, 3\"),\n\n    @(13, 1, \"46\u00f73=15, 1\"),\n    @(13, 2, \"87\u00f72=43, 1\"),\n    @(13, 3, \"87\u00f74=21, 3\"),\n    @(13, 4, \"46\u00f77=6, 4\"),\n    @(13, 5, \"78\u00f72=39, 0\"),\n\n    @(17, 1, \"50\u00f78=6, 2\"),\n    @(17, 2, \"89\u00f76=14, 5\"),\n    @(17, 3, \"46\u00f79=5, 1\"),\n    @(17, 4, \"36\u00f77=5, 1\"),\n    @(17, 5, \"57\u00f72=28, 1\")\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $text = $u[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
